# Update the workbook for the 2022-07-11 data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-07-11"

# Update the header label in I1 to match the new "through" date.
$ws.Range("I1").Value = "2022 (through 07-11)"

# Update the July figure (row 8) for the "Total"/current-year column (I).
$ws.Range("I8").Value = 63

# Update the yearly Total row (row 14) for column I.
$ws.Range("I14").Value = 868
